$d = $word.ActiveDocument

# Locate the misspelled "Marcini" and fix it to "Marchini" by inserting
# the missing "h" after "Marc".
$rng = $d.Content
$rng.Find.Execute("Marcini", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$start = $rng.Start

$insPos = $start + 4
$insRange = $d.Range($insPos, $insPos)
$insRange.Text = "h"

# Word moves its "_GoBack" (last-edit-location) bookmark to the spot that
# was just edited; re-point it here (it currently sits before "haplotype
# effects drive the results" from the previous edit).
$d.Bookmarks.Add("_GoBack", $insRange)

# Collapse the bookmark to an insertion point immediately after the "h"
# that was just typed.
$afterH = $start + 5
$bmRange = $d.Range($afterH, $afterH)
$d.Bookmarks.Add("_GoBack", $bmRange)
